$p = $ppt.ActivePresentation
$s = $p.Slides.Add(7, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Diapo del profe"
$s.Shapes.Item(1).TextFrame.TextRange.LanguageID = "es-ES"
